$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values
$ws.Range("D2").Value = "'28.590.91"
$ws.Range("D3").Value = "'1.913.13"
$ws.Range("D5").Value = "'315.78"
$ws.Range("D6").Value = "'1.001"
$ws.Range("D7").Value = "'0.5125"
$ws.Range("D8").Value = "'0.3964"
$ws.Range("D9").Value = "'0.09771"
$ws.Range("D11").Value = "'42.18"
$ws.Range("D12").Value = "'6.552"
$ws.Range("D13").Value = "'21.24"
$ws.Range("D14").Value = "'1.918.13"
$ws.Range("D15").Value = "'7.593"
$ws.Range("D16").Value = "'1.001"
$ws.Range("D17").Value = "'0.00001141"
$ws.Range("D18").Value = "'94.00"
$ws.Range("D19").Value = "'0.06666"
$ws.Range("D20").Value = "'18.20"
$ws.Range("D22").Value = "'6.326"
$ws.Range("D23").Value = "'28.638.94"
$ws.Range("D24").Value = "'11.48"
$ws.Range("D25").Value = "'2.288"
$ws.Range("D26").Value = "'2.729"
$ws.Range("D27").Value = "'2.135.76"
$ws.Range("D28").Value = "'21.34"
$ws.Range("D29").Value = "'159.64"
$ws.Range("D30").Value = "'128.87"
$ws.Range("D32").Value = "'0.1081"
$ws.Range("D33").Value = "'5.753"
$ws.Range("D34").Value = "'3.646"
$ws.Range("D35").Value = "'9.905"
$ws.Range("D36").Value = "'0.06818"
$ws.Range("D37").Value = "'0.02445"
$ws.Range("D38").Value = "'1.272"
$ws.Range("D39").Value = "'0.2237"
$ws.Range("D40").Value = "'11.96"
$ws.Range("D41").Value = "'5.112"
$ws.Range("D42").Value = "'0.6460"
$ws.Range("D43").Value = "'1.195"
$ws.Range("D45").Value = "'13.75"
$ws.Range("D46").Value = "'0.6106"
$ws.Range("D49").Value = "'2.046"
$ws.Range("D50").Value = "'125.44"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("E3").Value = "  +5.59%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  +2.74%  "
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E11").Value = "  +3.27%  "
$ws.Range("E12").Value = "  +2.33%  "
$ws.Range("E13").Value = "  +4.09%  "
$ws.Range("E14").Value = "  +5.84%  "
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("E20").Value = "  +6.12%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  +7.27%  "
$ws.Range("E23").Value = "  +2.30%  "
$ws.Range("E24").Value = "  +3.83%  "
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("E26").Value = "  +14.50%  "
$ws.Range("E27").Value = "  +5.69%  "
$ws.Range("E28").Value = "  +4.01%  "
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("E31").Value = "  +7.25%  "
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("E33").Value = "  +3.60%  "
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("E35").Value = "  +11.96%  "
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("E37").Value = "  +5.26%  "
$ws.Range("E38").Value = "  +8.75%  "
$ws.Range("E39").Value = "  +4.75%  "
$ws.Range("E40").Value = "  +6.27%  "
$ws.Range("E41").Value = "  +3.76%  "
$ws.Range("E42").Value = "  +4.91%  "
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  +5.06%  "
$ws.Range("E47").Value = "  +3.26%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  +5.91%  "
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("E51").Value = "  +3.23%  "
